$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 461.7143
$ws.Range("I4").Value = 461.7143
$ws.Range("K4").Value = 461.7143
$ws.Range("M4").Value = -347.7143
$ws.Range("H38").Value = 76285.586
$ws.Range("J38").Value = 91499.89999999999
$ws.Range("L38").Value = 274499.7
$ws.Range("N38").Value = -275243.7
$ws.Range("H41").Value = 1490.8
$ws.Range("I41").Value = 329.33334
$ws.Range("K41").Value = 329.33334
$ws.Range("M41").Value = 110.66666
$ws.Range("H58").Value = 640.2857
$ws.Range("I58").Value = 640.2857
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 1920.8571
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1770.8571
$ws.Range("N58").ClearContents()
$ws.Range("H62").Value = 6116.8945
$ws.Range("I62").Value = 6130.5884
$ws.Range("K62").Value = 6130.5884
$ws.Range("M62").Value = -5506.5884
$ws.Range("H65").Value = 6116.8945
$ws.Range("I65").Value = 6130.5884
$ws.Range("K65").Value = 30652.942
$ws.Range("M65").Value = -27532.942
$ws.Range("H76").Value = 4463
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 4463
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H86").Value = 400001380
$ws.Range("I86").Value = 1000000000
$ws.Range("J86").Value = 2282.3333
$ws.Range("K86").Value = 1000000000
$ws.Range("L86").Value = 2282.3333
$ws.Range("M86").Value = -999998877
$ws.Range("N86").Value = -4528.3333
$ws.Range("H88").Value = 25250700
$ws.Range("I88").Value = 50500000
$ws.Range("J88").Value = 1400
$ws.Range("K88").Value = 50500000
$ws.Range("L88").Value = 1400
$ws.Range("M88").Value = -50499594
$ws.Range("N88").Value = -2212
$ws.Range("H89").Value = 400001380
$ws.Range("I89").Value = 1000000000
$ws.Range("J89").Value = 2282.3333
$ws.Range("K89").Value = 5000000000
$ws.Range("L89").Value = 11411.6665
$ws.Range("M89").Value = -4999994384
$ws.Range("N89").Value = -22643.6665
$ws.Range("H91").Value = 25250700
$ws.Range("I91").Value = 50500000
$ws.Range("J91").Value = 1400
$ws.Range("K91").Value = 50500000
$ws.Range("L91").Value = 1400
$ws.Range("M91").Value = -50498596
$ws.Range("N91").Value = -4208
$ws.Range("H98").Value = 3006.8948
$ws.Range("I98").Value = 3239.7878
$ws.Range("J98").Value = 1469.8
$ws.Range("K98").Value = 3239.7878
$ws.Range("L98").Value = 1469.8
$ws.Range("M98").Value = -1741.7878
$ws.Range("N98").Value = -4465.8
$ws.Range("H122").Value = 3006.8948
$ws.Range("I122").Value = 3239.7878
$ws.Range("J122").Value = 1469.8
$ws.Range("K122").Value = 9719.3634
$ws.Range("L122").Value = 4409.4
$ws.Range("M122").Value = -7269.3634
$ws.Range("N122").Value = -9309.4
$ws.Range("H133").Value = 105992.5
$ws.Range("J133").Value = 105992.5
$ws.Range("L133").Value = 105992.5
$ws.Range("N133").Value = -116112.5
$ws.Range("H137").Value = 2347.257
$ws.Range("I137").Value = 2261.1875
$ws.Range("J137").Value = 2419.7368
$ws.Range("K137").Value = 6783.5625
$ws.Range("L137").Value = 7259.2104
$ws.Range("M137").Value = -4233.5625
$ws.Range("N137").Value = -12359.2104
$ws.Range("H138").Value = 372876.56
$ws.Range("I138").Value = 4047.2222
$ws.Range("J138").Value = 429138.66
$ws.Range("K138").Value = 12141.6666
$ws.Range("L138").Value = 1287415.98
$ws.Range("M138").Value = -7001.6666
$ws.Range("N138").Value = -1297695.98

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11055.923
$ws.Range("I32").Value = 7726.0625
$ws.Range("J32").Value = 35273.09
$ws.Range("K32").Value = 7726.0625
$ws.Range("L32").Value = 35273.09
$ws.Range("M32").Value = -7439.0625
$ws.Range("N32").Value = -35847.09
$ws.Range("H45").Value = 15272.667
$ws.Range("I45").Value = 21153.305
$ws.Range("K45").Value = 21153.305
$ws.Range("M45").Value = -20776.305
$ws.Range("H60").Value = 21037
$ws.Range("I60").Value = 21037
$ws.Range("K60").Value = 21037
$ws.Range("M60").Value = -20304
$ws.Range("H74").Value = 1117377.4
$ws.Range("I74").Value = 1859629.4
$ws.Range("J74").Value = 3999.5
$ws.Range("K74").Value = 1859629.4
$ws.Range("L74").Value = 3999.5
$ws.Range("M74").Value = -1858755.4
$ws.Range("N74").Value = -5747.5
$ws.Range("H77").Value = 1117377.4
$ws.Range("I77").Value = 1859629.4
$ws.Range("J77").Value = 3999.5
$ws.Range("K77").Value = 9298147
$ws.Range("L77").Value = 19997.5
$ws.Range("M77").Value = -9293779
$ws.Range("N77").Value = -28733.5
$ws.Range("H108").Value = 49898
$ws.Range("J108").Value = 49898
$ws.Range("L108").Value = 49898
$ws.Range("N108").Value = -57578
$ws.Range("H132").Value = 2673.879
$ws.Range("I132").Value = 2064.5
$ws.Range("J132").Value = 3611.3845
$ws.Range("K132").Value = 6193.5
$ws.Range("L132").Value = 10834.1535
$ws.Range("M132").Value = -3663.5
$ws.Range("N132").Value = -15894.1535

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1491.5454
$ws.Range("I107").Value = 1369.2858
$ws.Range("K107").Value = 1369.2858
$ws.Range("M107").Value = 550.7141999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 14297479
$ws.Range("I62").Value = 33338266
$ws.Range("J62").Value = 16887.5
$ws.Range("K62").Value = 33338266
$ws.Range("L62").Value = 16887.5
$ws.Range("M62").Value = -33337642
$ws.Range("N62").Value = -18135.5
$ws.Range("H65").Value = 14297479
$ws.Range("I65").Value = 33338266
$ws.Range("J65").Value = 16887.5
$ws.Range("K65").Value = 166691330
$ws.Range("L65").Value = 84437.5
$ws.Range("M65").Value = -166688210
$ws.Range("N65").Value = -90677.5
$ws.Range("H132").Value = 12199252
$ws.Range("I132").Value = 13516793
$ws.Range("K132").Value = 40550379
$ws.Range("M132").Value = -40547849
$ws.Range("H134").Value = 3025.353
$ws.Range("I134").Value = 2218.3572
$ws.Range("J134").Value = 6791.3335
$ws.Range("K134").Value = 6655.071599999999
$ws.Range("L134").Value = 20374.0005
$ws.Range("M134").Value = -4120.071599999999
$ws.Range("N134").Value = -25444.0005
$ws.Range("H141").Value = 594004.4
$ws.Range("J141").Value = 594004.4
$ws.Range("L141").Value = 594004.4
$ws.Range("N141").Value = -604364.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1195
$ws.Range("J107").Value = 1318.4615
$ws.Range("L107").Value = 3955.3845
$ws.Range("N107").Value = -7795.3845
$ws.Range("H122").Value = 1874.9166
$ws.Range("I122").Value = 1662.875
$ws.Range("J122").Value = 1980.9375
$ws.Range("K122").Value = 14965.875
$ws.Range("L122").Value = 17828.4375
$ws.Range("M122").Value = -12515.875
$ws.Range("N122").Value = -22728.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 3400
$ws.Range("J7").Value = 3400
$ws.Range("L7").Value = 3400
$ws.Range("N7").Value = -3624
$ws.Range("H8").Value = 3400
$ws.Range("J8").Value = 3400
$ws.Range("L8").Value = 3400
$ws.Range("N8").Value = -3678
$ws.Range("H28").Value = 25000
$ws.Range("J28").Value = 25000
$ws.Range("L28").Value = 25000
$ws.Range("N28").Value = -25384
$ws.Range("H97").Value = 823
$ws.Range("J97").Value = 799
$ws.Range("L97").Value = 799
$ws.Range("N97").Value = -1791
$ws.Range("H132").Value = 2102.5
$ws.Range("I132").Value = 1835.3334
$ws.Range("J132").Value = 2503.25
$ws.Range("K132").Value = 5506.0002
$ws.Range("L132").Value = 7509.75
$ws.Range("M132").Value = -2976.0002
$ws.Range("N132").Value = -12569.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H123").Value = 66189.25
$ws.Range("J123").Value = 66189.25
$ws.Range("L123").Value = 66189.25
$ws.Range("N123").Value = -75989.25
$ws.Range("H132").Value = 3644.0588
$ws.Range("I132").Value = 3282.0715
$ws.Range("K132").Value = 9846.2145
$ws.Range("M132").Value = -7316.2145
$ws.Range("H136").Value = 15604
$ws.Range("I136").Value = 11906
$ws.Range("J136").Value = 23000
$ws.Range("K136").Value = 35718
$ws.Range("L136").Value = 69000
$ws.Range("M136").Value = -33168
$ws.Range("N136").Value = -74100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7001.5
$ws.Range("J62").Value = 7001.5
$ws.Range("L62").Value = 7001.5
$ws.Range("N62").Value = -8249.5
$ws.Range("H65").Value = 7001.5
$ws.Range("J65").Value = 7001.5
$ws.Range("L65").Value = 35007.5
$ws.Range("N65").Value = -41247.5
$ws.Range("H81").Value = 6716.6665
$ws.Range("I81").Value = 5783
$ws.Range("J81").Value = 7650.3335
$ws.Range("K81").Value = 11566
$ws.Range("L81").Value = 15300.667
$ws.Range("M81").Value = -10505
$ws.Range("N81").Value = -17422.667
$ws.Range("H84").Value = 6716.6665
$ws.Range("I84").Value = 5783
$ws.Range("J84").Value = 7650.3335
$ws.Range("K84").Value = 57830
$ws.Range("L84").Value = 76503.33499999999
$ws.Range("M84").Value = -52526
$ws.Range("N84").Value = -87111.33499999999
$ws.Range("H132").Value = 2774.97
$ws.Range("I132").Value = 2762.8555
$ws.Range("J132").Value = 2834.1177
$ws.Range("K132").Value = 8288.566500000001
$ws.Range("L132").Value = 8502.3531
$ws.Range("M132").Value = -5758.566500000001
$ws.Range("N132").Value = -13562.3531
$ws.Range("H136").Value = 2250.9033
$ws.Range("I136").Value = 2011.75
$ws.Range("J136").Value = 4483
$ws.Range("K136").Value = 6035.25
$ws.Range("L136").Value = 13449
$ws.Range("M136").Value = -3485.25
$ws.Range("N136").Value = -18549
